# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit calculation sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 30400.666
$ws.Range("I9").Value = 45500
$ws.Range("J9").Value = 202
$ws.Range("K9").Value = 45500
$ws.Range("L9").Value = 202
$ws.Range("M9").Value = -45331
$ws.Range("N9").Value = -540

$ws.Range("H17").Value = 742921.3
$ws.Range("J17").Value = 742921.3
$ws.Range("L17").Value = 2228763.9
$ws.Range("N17").Value = -2229099.9

$ws.Range("H86").Value = 6996
$ws.Range("I86").Value = 9994
$ws.Range("J86").Value = 5497
$ws.Range("K86").Value = 9994
$ws.Range("L86").Value = 5497
$ws.Range("M86").Value = -8871
$ws.Range("N86").Value = -7743

$ws.Range("H88").Value = 1249.091
$ws.Range("J88").Value = 1274.2
$ws.Range("L88").Value = 1274.2
$ws.Range("N88").Value = -2086.2

$ws.Range("H89").Value = 6996
$ws.Range("I89").Value = 9994
$ws.Range("J89").Value = 5497
$ws.Range("K89").Value = 49970
$ws.Range("L89").Value = 27485
$ws.Range("M89").Value = -44354
$ws.Range("N89").Value = -38717

$ws.Range("H91").Value = 1249.091
$ws.Range("J91").Value = 1274.2
$ws.Range("L91").Value = 1274.2
$ws.Range("N91").Value = -4082.2

$ws.Range("H96").Value = 2135.375
$ws.Range("I96").Value = 2220.6667
$ws.Range("J96").Value = 1879.5
$ws.Range("K96").Value = 6662.000100000001
$ws.Range("L96").Value = 5638.5
$ws.Range("M96").Value = -5289.000100000001
$ws.Range("N96").Value = -8384.5

$ws.Range("H97").Value = 21239
$ws.Range("J97").Value = 28010.334
$ws.Range("L97").Value = 84031.00199999999
$ws.Range("N97").Value = -85023.00199999999

$ws.Range("H111").Value = 1121.9333
$ws.Range("I111").Value = 711.7273
$ws.Range("J111").Value = 2250
$ws.Range("K111").Value = 2135.1819
$ws.Range("L111").Value = 6750
$ws.Range("M111").Value = 931.8181
$ws.Range("N111").Value = -12884

$ws.Range("H116").Value = 50459.92
$ws.Range("J116").Value = 9999.833000000001
$ws.Range("L116").Value = 9999.833000000001
$ws.Range("N116").Value = -16883.833

$ws.Range("H135").Value = 1351.875
$ws.Range("I135").Value = 993
$ws.Range("J135").Value = 1710.75
$ws.Range("K135").Value = 8937
$ws.Range("L135").Value = 15396.75
$ws.Range("M135").Value = -6402
$ws.Range("N135").Value = -20466.75

$ws.Range("H137").Value = 9619.333000000001
$ws.Range("I137").Value = 1062.5
$ws.Range("J137").Value = 13897.75
$ws.Range("K137").Value = 3187.5
$ws.Range("L137").Value = 41693.25
$ws.Range("M137").Value = -637.5
$ws.Range("N137").Value = -46793.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4586.125
$ws.Range("I2").Value = 4586.125
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4586.125
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -4473.125

$ws.Range("H55").Value = 12889.4
$ws.Range("I55").Value = 8149.6665
$ws.Range("K55").Value = 8149.6665
$ws.Range("M55").Value = -7834.6665

$ws.Range("H88").Value = 1502
$ws.Range("I88").Value = 920.8333
$ws.Range("J88").Value = 1850.7
$ws.Range("K88").Value = 920.8333
$ws.Range("L88").Value = 1850.7
$ws.Range("M88").Value = -514.8333
$ws.Range("N88").Value = -2662.7

$ws.Range("H91").Value = 1502
$ws.Range("I91").Value = 920.8333
$ws.Range("J91").Value = 1850.7
$ws.Range("K91").Value = 920.8333
$ws.Range("L91").Value = 1850.7
$ws.Range("M91").Value = 483.1667
$ws.Range("N91").Value = -4658.7

$ws.Range("H97").Value = 583.4
$ws.Range("I97").Value = 553.7778
$ws.Range("K97").Value = 553.7778
$ws.Range("M97").Value = -57.77779999999996

$ws.Range("H102").Value = 5312.2
$ws.Range("I102").Value = 4516.8184
$ws.Range("K102").Value = 4516.8184
$ws.Range("M102").Value = -2894.8184

$ws.Range("H116").Value = 4586.125
$ws.Range("I116").Value = 4586.125
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4586.125
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -2292.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4586.125
$ws.Range("I3").Value = 4586.125
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4586.125
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -4472.125

$ws.Range("H22").Value = 426501.94
$ws.Range("I22").Value = 748.2308
$ws.Range("K22").Value = 748.2308
$ws.Range("M22").Value = -575.2308

$ws.Range("H94").Value = 1343.8572
$ws.Range("I94").Value = 1827
$ws.Range("K94").Value = 1827
$ws.Range("M94").Value = -1376

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4332.6665
$ws.Range("J99").Value = 7241
$ws.Range("L99").Value = 7241
$ws.Range("N99").Value = -10237

$ws.Range("H126").Value = 4332.6665
$ws.Range("J126").Value = 7241
$ws.Range("L126").Value = 21723
$ws.Range("N126").Value = -26663

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 219.16667
$ws.Range("I47").Value = 163.33333
$ws.Range("J47").Value = 275
$ws.Range("K47").Value = 489.99999
$ws.Range("L47").Value = 825
$ws.Range("M47").Value = -58.99998999999997
$ws.Range("N47").Value = -1687

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 646
$ws.Range("I107").Value = 769.875
$ws.Range("J107").Value = 425.77777
$ws.Range("K107").Value = 769.875
$ws.Range("L107").Value = 425.77777
$ws.Range("M107").Value = 1150.125
$ws.Range("N107").Value = -4265.77777

$ws.Range("H132").Value = 2666.054
$ws.Range("I132").Value = 2460.3235
$ws.Range("K132").Value = 7380.970499999999
$ws.Range("M132").Value = -4850.970499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 11364.685
$ws.Range("J43").Value = 17999.8
$ws.Range("L43").Value = 17999.8
$ws.Range("N43").Value = -18385.8

$ws.Range("H87").Value = 97303.42999999999
$ws.Range("J87").Value = 97303.42999999999
$ws.Range("L87").Value = 97303.42999999999
$ws.Range("N87").Value = -99549.42999999999

$ws.Range("H90").Value = 97303.42999999999
$ws.Range("J90").Value = 97303.42999999999
$ws.Range("L90").Value = 291910.29
$ws.Range("N90").Value = -303142.29

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 2499.8
$ws.Range("I10").Value = 2166.6667
$ws.Range("J10").Value = 2999.5
$ws.Range("K10").Value = 2166.6667
$ws.Range("L10").Value = 2999.5
$ws.Range("M10").Value = -1997.6667
$ws.Range("N10").Value = -3337.5

$ws.Range("H80").Value = 35198.8
$ws.Range("J80").Value = 40748.5
$ws.Range("L80").Value = 40748.5
$ws.Range("N80").Value = -42744.5

$ws.Range("H81").Value = 9200.409
$ws.Range("I81").Value = 27352.5
$ws.Range("J81").Value = 5166.6113
$ws.Range("K81").Value = 54705
$ws.Range("L81").Value = 10333.2226
$ws.Range("M81").Value = -53644
$ws.Range("N81").Value = -12455.2226

$ws.Range("H83").Value = 35198.8
$ws.Range("J83").Value = 40748.5
$ws.Range("L83").Value = 122245.5
$ws.Range("N83").Value = -132229.5

$ws.Range("H84").Value = 9200.409
$ws.Range("I84").Value = 27352.5
$ws.Range("J84").Value = 5166.6113
$ws.Range("K84").Value = 273525
$ws.Range("L84").Value = 51666.113
$ws.Range("M84").Value = -268221
$ws.Range("N84").Value = -62274.113

$ws.Range("H96").Value = 2893.875
$ws.Range("J96").Value = 3146.3333
$ws.Range("L96").Value = 3146.3333
$ws.Range("N96").Value = -5892.3333

$ws.Range("H113").Value = 1152.5834
$ws.Range("I113").Value = 840.2143
$ws.Range("K113").Value = 2520.6429
$ws.Range("M113").Value = -350.6428999999998

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0
